$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1560
$ws.Range("C3").Value = "Chocolate"
$ws.Range("D3").Value = 890
$ws.Range("C4").Value = "Chocolate"
$ws.Range("D4").Value = 1780
$ws.Range("C5").Value = "Chicle"
$ws.Range("D5").Value = 160
$ws.Range("C6").Value = "Chicle"
$ws.Range("D6").Value = 240
$ws.Range("C7").Value = "Caramelo"
$ws.Range("D7").Value = 150
$ws.Range("C8").Value = "Chocolate"
$ws.Range("D8").Value = 890
$ws.Range("C9").Value = "Galletitas"
$ws.Range("D9").Value = 2080
$ws.Range("C10").Value = "Caramelo"
$ws.Range("D10").Value = 100
$ws.Range("E11").Value = "QR"
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = "QR"
$ws.Range("C13").Value = "Chocolate"
$ws.Range("D13").Value = 1780
$ws.Range("E13").Value = "Debito"
$ws.Range("D14").Value = 2250
$ws.Range("E14").Value = "QR"
$ws.Range("D15").Value = 4450
$ws.Range("E15").Value = "QR"
$ws.Range("C16").Value = "Galletitas"
$ws.Range("D16").Value = 1040
$ws.Range("E16").Value = "QR"
$ws.Range("D17").Value = 320
$ws.Range("E17").Value = "Debito"
$ws.Range("C18").Value = "Chocolate"
$ws.Range("D18").Value = 3560
$ws.Range("E18").Value = "QR"
$ws.Range("C19").Value = "Chicle"
$ws.Range("D19").Value = 240
$ws.Range("E19").Value = "Debito"
